$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '22.425.90'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.15%  '
# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.573.34'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.07%  '
# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.003'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.07%  '
# Row 5
$ws.Range('E5').Value = '  -0.04%  '
# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '291.28'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.13%  '
# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3731'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.81%  '
# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '49.97'
$ws.Range('D8').Style = 'Normal'
# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3397'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.85%  '
# Row 10
$ws.Range('B10').Value = 'Dogecoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07574'
$ws.Range('D10').Style = 'Normal'
# Row 11
$ws.Range('B11').Value = 'Polygon'
$ws.Range('C11').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.144'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.90%  '
# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.003'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.00%  '
# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '21.34'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.64%  '
# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.013'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.19%  '
# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.958'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.23%  '
# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.572.05'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.05%  '
# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001124'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.84%  '
# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '90.94'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.98%  '
# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06751'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.07%  '
# Row 20
$ws.Range('E20').Value = '  +0.08%  '
# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.305'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.69%  '
# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '16.33'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.91%  '
# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '12.16'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.34%  '
# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '22.420.40'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.10%  '
# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.349'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.99%  '
# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.682'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.41%  '
# Row 27
$ws.Range('E27').Value = '  -0.37%  '
# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '148.58'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.86%  '
# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.011'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.58%  '
# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '125.65'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.44%  '
# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.747.94'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.05%  '
# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.050'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +6.58%  '
# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.194'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.60%  '
# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.980'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.35%  '
# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '9.816'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.37%  '
# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.08401'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.47%  '
# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.377'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +3.09%  '
# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02483'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.45%  '
# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.2292'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.07%  '
# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.06530'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.56%  '
# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.473'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.97%  '
# Row 42
$ws.Range('E42').Value = '  -1.05%  '
# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.6238'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.57%  '
# Row 44
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '14.13'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.55%  '
# Row 45
$ws.Range('B45').Value = 'Frax'
$ws.Range('C45').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.001'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.05%  '
# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.815'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.83%  '
# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5815'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.90%  '
# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '129.37'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +3.00%  '
# Row 49
$ws.Range('E49').Value = '  -0.48%  '
# Row 50
$ws.Range('E50').Value = '  -6.31%  '
# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.07334'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.10%  '
